$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# B10/C10/D10 hold the MSME / SMEs / MSMEs "Enterprises (% of total)" figures
# as text (they look numeric but are stored as plain strings in the sheet).
# Prefixing with an apostrophe forces the new values to stay text instead of
# being auto-converted to numbers, matching the original cell type.
$ws.Range("B10").Value = "'88.06"
$ws.Range("C10").Value = "'11.69"
$ws.Range("D10").Value = "'99.75"

# Re-apply the default "Normal" cell style so the cells keep their original
# (unformatted) appearance instead of picking up the quote-prefix style that
# Excel creates when forcing text entry.
$ws.Range("B10:D10").Style = "Normal"
